$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "VaR YYDS"
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").Font.Color = 0

$ws.Range("B5").Value = "VaR calculator for Stock Portfolio through Monte Carlo and Geometric Brownian Motion"
$ws.Range("F22").Value = "https://github.com/fayeyeye"

$ws.Application.ActiveWindow.Zoom = 175
$ws.Range("B9").Select()
